$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.398.04"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "1.880.14"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'0.7156"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").Value = "'243.58"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Value = "'0.07948"
$ws.Range("E8").Value = "  +1.90%  "

$ws.Range("D9").Value = "'0.3144"
$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("D10").Value = "'24.94"
$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("D11").Value = "'0.08154"
$ws.Range("E11").Value = "  -2.70%  "

$ws.Range("D12").Value = "1.906.06"
$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").Value = "'5.246"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").Value = "'94.66"
$ws.Range("E14").Value = "  +3.83%  "

$ws.Range("D15").Value = "'0.7088"
$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("E16").Value = "  +4.54%  "

$ws.Range("D17").Value = "'0.000008421"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "29.418.72"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").Value = "'253.31"
$ws.Range("E19").Value = "  +5.46%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.143.29"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'13.31"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "'7.741"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").Value = "'0.1585"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").Value = "'9.070"
$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("D27").Value = "'162.53"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").Value = "'18.93"
$ws.Range("E28").Value = "  +2.36%  "

$ws.Range("D29").Value = "'1.505"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").Value = "'4.412"
$ws.Range("E30").Value = "  -0.27%  "

$ws.Range("D31").Value = "'4.284"
$ws.Range("E31").Value = "  -1.39%  "

$ws.Range("D32").Value = "'1.232"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").Value = "'0.05324"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").Value = "'1.949"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").Value = "'0.7575"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("D36").Value = "'1.178"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").Value = "'2.697"
$ws.Range("E37").Value = "  +0.50%  "

$ws.Range("D38").Value = "'0.01896"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").Value = "1.276.01"
$ws.Range("E39").Value = "  +3.09%  "

$ws.Range("D40").Value = "'2.765"
$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("D41").Value = "'6.411"
$ws.Range("E41").Value = "  -1.93%  "

$ws.Range("D42").Value = "'112.72"
$ws.Range("E42").Value = "  +3.37%  "

$ws.Range("D43").Value = "'74.30"
$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("D44").Value = "'0.9036"
$ws.Range("E44").Value = "  +1.17%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000131"
$ws.Range("E45").Value = "  +2.46%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").Value = "2.037.64"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("D48").Value = "'1.805"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("D49").Value = "'0.5211"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").Value = "'9.507"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("D51").Value = "'0.4352"
$ws.Range("E51").Value = "  +0.11%  "
